# Generate Report for Archive
#
# 1) The localization status of the two handed-off files moves from
#    "Ready for handoff" to "In Translation" (Overview sheet columns E/F,
#    and the per-language "Status" column on the zh-cn / de-de sheets).
# 2) The Status/"zh-cn"/"de-de" columns are narrowed to better fit the new
#    (shorter) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
# Overview sheet: zh-cn status is column E, de-de status is column F, for
# both data rows (2 and 3).
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Per-language sheets: Status is column C, for both data rows.
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the status columns ------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
